$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.829945333333333
$ws.Range("H2").Value = 5.489835999999999
$ws.Range("I2").Value = 0.4190796720210465
$ws.Range("J2").Value = 0.4190796720210465
$ws.Range("M2").Value = 0.1994653333333334
$ws.Range("N2").Value = 0.598396
$ws.Range("O2").Value = 0.01676579960230272
$ws.Range("P2").Value = 0.01676579960230271
$ws.Range("Q2").Value = 0.3650106558951111
$ws.Range("R2").Value = 3.285095903056
$ws.Range("S2").Value = 0.007026205798503614
$ws.Range("T2").Value = 0.007026205798503612
$ws.Range("G3").Value = 1.829945333333333
$ws.Range("H3").Value = 5.489835999999999
$ws.Range("I3").Value = 0.4190796720210465
$ws.Range("J3").Value = 0.4190796720210465
$ws.Range("O3").Value = 0.03203779682023726
$ws.Range("P3").Value = 0.03203779682023726
$ws.Range("Q3").Value = 0.6974995233262221
$ws.Range("R3").Value = 6.277495709935999
$ws.Range("S3").Value = 0.01342638938370196
$ws.Range("T3").Value = 0.01342638938370196
$ws.Range("G4").Value = 1.829945333333333
$ws.Range("H4").Value = 5.489835999999999
$ws.Range("I4").Value = 0.4190796720210465
$ws.Range("J4").Value = 0.4190796720210465
$ws.Range("M4").Value = 0.2888043333333333
$ws.Range("N4").Value = 0.8664129999999999
$ws.Range("O4").Value = 0.02427507324719734
$ws.Range("P4").Value = 0.02427507324719734
$ws.Range("Q4").Value = 0.5284961420297776
$ws.Range("R4").Value = 4.756465278267999
$ws.Range("S4").Value = 0.01017318973472234
$ws.Range("T4").Value = 0.01017318973472234
$ws.Range("G5").Value = 1.829945333333333
$ws.Range("H5").Value = 5.489835999999999
$ws.Range("I5").Value = 0.4190796720210465
$ws.Range("J5").Value = 0.4190796720210465
$ws.Range("M5").Value = 11.02772766666667
$ws.Range("N5").Value = 33.083183
$ws.Range("O5").Value = 0.9269213303302627
$ws.Range("P5").Value = 0.9269213303302626
$ws.Range("Q5").Value = 20.18013878088755
$ws.Range("R5").Value = 181.621249027988
$ws.Range("S5").Value = 0.3884538871041185
$ws.Range("T5").Value = 0.3884538871041185
$ws.Range("I6").Value = 0.2833335737960661
$ws.Range("J6").Value = 0.2833335737960661
$ws.Range("M6").Value = 0.1994653333333334
$ws.Range("N6").Value = 0.598396
$ws.Range("O6").Value = 0.01676579960230272
$ws.Range("P6").Value = 0.01676579960230271
$ws.Range("Q6").Value = 0.2467783109346667
$ws.Range("R6").Value = 2.221004798412
$ws.Range("S6").Value = 0.004750313918869092
$ws.Range("T6").Value = 0.00475031391886909
$ws.Range("I7").Value = 0.2833335737960661
$ws.Range("J7").Value = 0.2833335737960661
$ws.Range("O7").Value = 0.03203779682023726
$ws.Range("P7").Value = 0.03203779682023726
$ws.Range("S7").Value = 0.009077383469630064
$ws.Range("T7").Value = 0.009077383469630064
$ws.Range("I8").Value = 0.2833335737960661
$ws.Range("J8").Value = 0.2833335737960661
$ws.Range("M8").Value = 0.2888043333333333
$ws.Range("N8").Value = 0.8664129999999999
$ws.Range("O8").Value = 0.02427507324719734
$ws.Range("P8").Value = 0.02427507324719734
$ws.Range("Q8").Value = 0.3573084323956666
$ws.Range("R8").Value = 3.215775891561
$ws.Range("S8").Value = 0.006877943257289697
$ws.Range("T8").Value = 0.006877943257289696
$ws.Range("I9").Value = 0.2833335737960661
$ws.Range("J9").Value = 0.2833335737960661
$ws.Range("M9").Value = 11.02772766666667
$ws.Range("N9").Value = 33.083183
$ws.Range("O9").Value = 0.9269213303302627
$ws.Range("P9").Value = 0.9269213303302626
$ws.Range("Q9").Value = 13.64349364147233
$ws.Range("R9").Value = 122.791442773251
$ws.Range("S9").Value = 0.2626279331502772
$ws.Range("T9").Value = 0.2626279331502772
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.1530633333333333
$ws.Range("H10").Value = 0.45919
$ws.Range("I10").Value = 0.03505335944376924
$ws.Range("J10").Value = 0.03505335944376924
$ws.Range("M10").Value = 0.1994653333333334
$ws.Range("N10").Value = 0.598396
$ws.Range("O10").Value = 0.01676579960230272
$ws.Range("P10").Value = 0.01676579960230271
$ws.Range("Q10").Value = 0.03053082880444445
$ws.Range("R10").Value = 0.27477745924
$ws.Range("S10").Value = 0.0005876975998217204
$ws.Range("T10").Value = 0.0005876975998217203
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.1530633333333333
$ws.Range("H11").Value = 0.45919
$ws.Range("I11").Value = 0.03505335944376924
$ws.Range("J11").Value = 0.03505335944376924
$ws.Range("O11").Value = 0.03203779682023726
$ws.Range("P11").Value = 0.03203779682023726
$ws.Range("Q11").Value = 0.05834141604888889
$ws.Range("R11").Value = 0.52507274444
$ws.Range("S11").Value = 0.001123032407726224
$ws.Range("T11").Value = 0.001123032407726224
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.1530633333333333
$ws.Range("H12").Value = 0.45919
$ws.Range("I12").Value = 0.03505335944376924
$ws.Range("J12").Value = 0.03505335944376924
$ws.Range("M12").Value = 0.2888043333333333
$ws.Range("N12").Value = 0.8664129999999999
$ws.Range("O12").Value = 0.02427507324719734
$ws.Range("P12").Value = 0.02427507324719734
$ws.Range("Q12").Value = 0.0442053539411111
$ws.Range("R12").Value = 0.3978481854699999
$ws.Range("S12").Value = 0.0008509228680578349
$ws.Range("T12").Value = 0.0008509228680578346
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.1530633333333333
$ws.Range("H13").Value = 0.45919
$ws.Range("I13").Value = 0.03505335944376924
$ws.Range("J13").Value = 0.03505335944376924
$ws.Range("M13").Value = 11.02772766666667
$ws.Range("N13").Value = 33.083183
$ws.Range("O13").Value = 0.9269213303302627
$ws.Range("P13").Value = 0.9269213303302626
$ws.Range("Q13").Value = 1.687940755752222
$ws.Range("R13").Value = 15.19146680177
$ws.Range("S13").Value = 0.03249170656816346
$ws.Range("T13").Value = 0.03249170656816346
$ws.Range("G14").Value = 1.146373333333333
$ws.Range("H14").Value = 3.43912
$ws.Range("I14").Value = 0.2625333947391181
$ws.Range("J14").Value = 0.2625333947391181
$ws.Range("M14").Value = 0.1994653333333334
$ws.Range("N14").Value = 0.598396
$ws.Range("O14").Value = 0.01676579960230272
$ws.Range("P14").Value = 0.01676579960230271
$ws.Range("Q14").Value = 0.2286617390577778
$ws.Range("R14").Value = 2.05795565152
$ws.Range("S14").Value = 0.004401582285108289
$ws.Range("T14").Value = 0.004401582285108288
$ws.Range("G15").Value = 1.146373333333333
$ws.Range("H15").Value = 3.43912
$ws.Range("I15").Value = 0.2625333947391181
$ws.Range("J15").Value = 0.2625333947391181
$ws.Range("O15").Value = 0.03203779682023726
$ws.Range("P15").Value = 0.03203779682023726
$ws.Range("Q15").Value = 0.4369501312355555
$ws.Range("R15").Value = 3.93255118112
$ws.Range("S15").Value = 0.008410991559179012
$ws.Range("T15").Value = 0.008410991559179012
$ws.Range("G16").Value = 1.146373333333333
$ws.Range("H16").Value = 3.43912
$ws.Range("I16").Value = 0.2625333947391181
$ws.Range("J16").Value = 0.2625333947391181
$ws.Range("M16").Value = 0.2888043333333333
$ws.Range("N16").Value = 0.8664129999999999
$ws.Range("O16").Value = 0.02427507324719734
$ws.Range("P16").Value = 0.02427507324719734
$ws.Range("Q16").Value = 0.3310775862844443
$ws.Range("R16").Value = 2.97969827656
$ws.Range("S16").Value = 0.006373017387127465
$ws.Range("T16").Value = 0.006373017387127465
$ws.Range("G17").Value = 1.146373333333333
$ws.Range("H17").Value = 3.43912
$ws.Range("I17").Value = 0.2625333947391181
$ws.Range("J17").Value = 0.2625333947391181
$ws.Range("M17").Value = 11.02772766666667
$ws.Range("N17").Value = 33.083183
$ws.Range("O17").Value = 0.9269213303302627
$ws.Range("P17").Value = 0.9269213303302626
$ws.Range("Q17").Value = 12.64189292432889
$ws.Range("R17").Value = 113.77703631896
$ws.Range("S17").Value = 0.2433478035077034
$ws.Range("T17").Value = 0.2433478035077033
